$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns C:G (all 7 rows), shifting nothing else - leaves only A and B
$ws.Range("C1:G7").Delete()

# Row 1: B1 keeps its bold/border style, gets the AAA080527HB7 tax id (previously in F1)
$ws.Range("B1").Value = 'AAA080527HB7'

# Row 2
$ws.Range("A2").Value = "'" + '2014'
$ws.Range("B2").Value = '{''Mision'': ''Darles un hogar digno en el que se encuentren felices y amados'', ''Vision'': ''Concluir la obra para poder trasladarlas, y poderles dar mejores espacios y areas verdes para que tengan mas comodidad, y acoger a muchos mas ya que la sociedad lo necesita, ya que el actual es muy pequeño ya que es una casa prestada'', ''Pagina web'': ''N/A'', ''Anio de autorizacion'': ''2013'', ''Socios o Asociados'': ''FELISA SANCHEZ FERNANDEZ, MARTIN HERNANDEZ, FRANCISCO LOPEZ AGUILAR, VICTOR ALFONSO GARCIA HERNANDEZ.'', ''Rubros autorizados'': [''A. Organizaciones Civiles y Fideicomisos Asistenciales.''], ''Plantilla laboral'': ''3'', ''Plantilla voluntariado'': ''6'', ''Monto total plantilla laboral'': ''30,000'', ''Activo'': ''3,998,689'', ''Pasivo'': ''0'', ''Capital'': ''3,998,689'', ''Gastos administracion'': ''2,999'', ''Gastos operacion'': ''432,572'', ''Gastos representacion'': ''0'', ''Ingresos'': [{''Ingresos'': [{''Especie '': ''ARTICULOS DE LIMPIEZA'', ''Efectivo '': '''', ''Tipo de donante'': ''Persona Moral Pública Nacional''}, {''Especie '': ''MEDICAMENTOS Y ANALISIS CLINICOS'', ''Efectivo '': '''', ''Tipo de donante'': ''Persona Moral Pública Nacional''}, {''Especie '': ''MUEBLES Y ENSERES MENORES'', ''Efectivo '': '''', ''Tipo de donante'': ''Persona Moral Pública Nacional''}, {''Especie '': ''ARTICULOS DE LIMPIEZA'', ''Efectivo '': '''', ''Tipo de donante'': ''Persona Moral Pública Nacional''}, {''Especie '': ''MEDICAMENTOS Y ANALISIS CLINICOS'', ''Efectivo '': '''', ''Tipo de donante'': ''Persona Moral Pública Nacional''}, {''Especie '': ''MUEBLES Y ENSERES MENORES'', ''Efectivo '': '''', ''Tipo de donante'': ''Persona Moral Pública Nacional''}]}, {''OtrosIngresos'': [{''Concepto'': ''REEMBOLSO POR CAJA'', ''Monto'': ''1,227''}, {''Concepto'': ''RIFAS'', ''Monto'': ''0''}, {''Concepto'': ''CUOTAS DE RECUPERACION'', ''Monto'': ''47,526''}, {''Concepto'': ''BAZAR'', ''Monto'': ''0''}, {''Concepto'': ''INTERESES BANCARIOS'', ''Monto'': ''236''}]}], ''Egresos'': [{''Egresos'': [{''Cargo ó Nombramiento'': ''NO SE OTORGARON PERCEPCIONES'', ''Total de percepciones netas'': ''0''}]}, {''MontosConceptos'': [{''Monto'': ''2,999'', ''Conceptos Ejercidos'': ''GASTOS FINANCIEROS''}, {''Monto'': ''5,478'', ''Conceptos Ejercidos'': ''OTROS''}, {''Monto'': ''2,159'', ''Conceptos Ejercidos'': ''SERVICIO TELEFONICO''}, {''Monto'': ''117,512'', ''Conceptos Ejercidos'': ''ACTIVOS FIJOS MOBILIARIO Y EQUIPO''}, {''Monto'': ''214'', ''Conceptos Ejercidos'': ''FLETES Y ACARREOS''}, {''Monto'': ''4,910'', ''Conceptos Ejercidos'': ''DIVERSOS''}, {''Monto'': ''0'', ''Conceptos Ejercidos'': ''PAPELERIA''}, {''Monto'': ''1,377'', ''Conceptos Ejercidos'': ''PAPELERIA Y UTILES''}, {''Monto'': ''2,467'', ''Conceptos Ejercidos'': ''ENERGIA ELECTRICA''}, {''Monto'': ''3,665'', ''Conceptos Ejercidos'': ''GAS''}, {''Monto'': ''138,489'', ''Conceptos Ejercidos'': ''MANTENIMIENTO E INSTALACION''}, {''Monto'': ''34,598'', ''Conceptos Ejercidos'': ''DESPENSA''}, {''Monto'': ''1,581,742'', ''Conceptos Ejercidos'': ''GASTOS POR OBRA CONSTRUCCION''}, {''Monto'': ''197,674'', ''Conceptos Ejercidos'': ''IVA ACREDITABLE''}, {''Monto'': ''2,731'', ''Conceptos Ejercidos'': ''MEDICAMENTO Y BOTIQUIN''}, {''Monto'': ''7,807'', ''Conceptos Ejercidos'': ''ALIMENTOS Y COMIDAS''}, {''Monto'': ''6,870'', ''Conceptos Ejercidos'': ''MANTENIMIENTO Y EQUIPO DE COMPUTO''}, {''Monto'': ''30,000'', ''Conceptos Ejercidos'': ''SUELDOS A ASIMILADOS''}]}, {''Donativos'': [{''Donativos'': ''0'', ''Beneficiarios'': ''NO SE OTORGARON''}]}], ''Actividades'': [{''Actividades'': []}]}'

# Row 3
$ws.Range("A3").Value = "'" + '2015'
$ws.Range("B3").Value = '{''Mision'': ''Mejorar la calidad de vida de los ancianos en condicion de abandono y maltrato de Guadalajara y la zona metropolitana brindando atencion integral al adulto mayor afreciendole un hogar digno con herramientas y personal necesario  para lograr un envejecimiento sano y activo'', ''Vision'': ''Ser una Organizacion reconocida a nivel nacional por excelente trato servicio y atencion que se le brinda a nuestros ancianos'', ''Pagina web'': ''http: asunciondemariaseniorshome.kerygs.com'', ''Anio de autorizacion'': ''2013'', ''Socios o Asociados'': ''FELISA  SANCHEZ FERNANDEZ / VICTOR MANUEL IBARRA GUERRERO/ ESTEBAN GARCIA REYES/JOSE ANTONIO SEPULVEDA/VICTOR ALFONSO GARCIA HERNANDEZ'', ''Rubros autorizados'': [''A. Organizaciones Civiles y Fideicomisos Asistenciales.''], ''Plantilla laboral'': ''3'', ''Plantilla voluntariado'': ''135'', ''Monto total plantilla laboral'': ''225,011'', ''Activo'': ''12,396,281'', ''Pasivo'': ''3,055'', ''Capital'': ''12,393,226'', ''Gastos administracion'': ''549,600'', ''Gastos operacion'': ''2,558,480'', ''Gastos representacion'': ''0'', ''Ingresos'': [{''Ingresos'': []}, {''OtrosIngresos'': [{''Concepto'': ''Cuotas de recuperacion'', ''Monto'': ''177,600''}, {''Concepto'': ''Colectas'', ''Monto'': ''20,000''}, {''Concepto'': ''Subsidios de Gobierno'', ''Monto'': ''1,500,000''}]}], ''Egresos'': [{''Egresos'': [{''Cargo ó Nombramiento'': ''.'', ''Total de percepciones netas'': ''0''}]}, {''MontosConceptos'': [{''Monto'': ''1,500,000'', ''Conceptos Ejercidos'': ''Construccion de Instalaciones''}]}, {''Donativos'': [{''Donativos'': ''0'', ''Beneficiarios'': ''180''}]}], ''Actividades'': [{''Actividades'': [{''Actividad'': ''Servicios Funerarios'', ''Número de Beneficiarios'': ''4'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''Tlaquepaque''}, {''Actividad'': ''estudios medicos especiales / terapias /analisis clinicos/homeopatia/ medicina alternativa'', ''Número de Beneficiarios'': ''180'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''Tlaquepaque''}, {''Actividad'': ''Casa Hogar al adulto mayor'', ''Número de Beneficiarios'': ''180'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''Tlaquepaque''}, {''Actividad'': ''Servicio de Comedor al adulto mayor'', ''Número de Beneficiarios'': ''180'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''Tlaquepaque''}, {''Actividad'': ''educacion basica'', ''Número de Beneficiarios'': ''2'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''Tlaquepaque''}, {''Actividad'': ''Talleres'', ''Número de Beneficiarios'': ''60'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''Tlaquepaque''}, {''Actividad'': ''atencion integral al adulto mayor'', ''Número de Beneficiarios'': ''180'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''Guadalajara''}]}]}'

# Row 4
$ws.Range("A4").Value = "'" + '2016'
$ws.Range("B4").Value = '{''Mision'': ''Mejorar la calidad de vida de los ancianos en condicion de abandono y maltrato de Guadalajara y la zona Metropolitana brindando atencion'', ''Vision'': ''ser una organizacion reconocida a nivel nacional por el excelente trato servicio y atencion que se les brinda a nuestros ancianos'', ''Pagina web'': ''http:asunciondemariaseniorshome.keyrgs.com'', ''Anio de autorizacion'': ''2013'', ''Socios o Asociados'': ''FELISA SANCHEZ FERNANDEZ PRESIDENTE DEL CONSEJO VICTOR ALFONSO GARCIA HERNANDEZ SECRETARIO ESTEBAN GARCIA REYES TESORERO JUANA RODRIGUEZ CONTRERAS VOCAL  MARIA LUISA VALDERRAMA CERVANTES VOCAL'', ''Rubros autorizados'': [''A. Organizaciones Civiles y Fideicomisos Asistenciales.''], ''Plantilla laboral'': ''13'', ''Plantilla voluntariado'': ''11'', ''Monto total plantilla laboral'': ''466,797'', ''Activo'': ''13,130,203'', ''Pasivo'': ''222,889'', ''Capital'': ''12,907,314'', ''Gastos administracion'': ''19,996'', ''Gastos operacion'': ''970,599'', ''Gastos representacion'': ''50,802'', ''Ingresos'': [{''Ingresos'': []}, {''OtrosIngresos'': [{''Concepto'': ''cuotas de recuperacion'', ''Monto'': ''1,019,776''}, {''Concepto'': ''subsidios de gobierno'', ''Monto'': ''399,994''}, {''Concepto'': ''colectas'', ''Monto'': ''41,000''}]}], ''Egresos'': [{''Egresos'': [{''Cargo ó Nombramiento'': ''ninguna'', ''Total de percepciones netas'': ''0''}]}, {''MontosConceptos'': [{''Monto'': ''1,929,952'', ''Conceptos Ejercidos'': ''asilo de ancianos''}]}, {''Donativos'': [{''Donativos'': ''0'', ''Beneficiarios'': ''Mercedes Sanchez Maria de laLuz Ramos Larios Luisa Macias Martinez Juan Jose Rito Villafan Vicente Guzman Galvan Antonio Castillo Castillo Patricia Araceli Coronado Hernandez Leoncio Gomez Jose Manuel Peredia Garcia Lorenzo Bautista Ramirez Maria Guadalupe Carrillo Soto Sara Rojas Trujillo Ramon Quintana Jose Miguel de Jesus Gonzalez Islas Marcelino Arellano Valin Vicente Ruiz Canel''}]}], ''Actividades'': [{''Actividades'': [{''Actividad'': ''atencion integral al adulto mayor'', ''Número de Beneficiarios'': ''79'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''34''}, {''Actividad'': ''servicio de comedor al adulto mayor'', ''Número de Beneficiarios'': ''79'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''tlaquepaque''}, {''Actividad'': ''estudios medicos especiales terapias analisis clinicos'', ''Número de Beneficiarios'': ''79'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''tlaquepaque''}, {''Actividad'': ''talleres'', ''Número de Beneficiarios'': ''60'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''tlaquepaque''}]}]}'

# Row 5
$ws.Range("A5").Value = "'" + '2017'
$ws.Range("B5").Value = '{''Mision'': ''Mejorar la calidad de vida de los ancianos en abandono y maltrato en Guadalajara y la zona metropolitana'', ''Vision'': ''Ser una organización reconocida a nivel nacional por el excelente trato y servicio y atención que se le brinda a los ancianos'', ''Pagina web'': ''http:asunciondemariaseniorshome.keyrgs.com'', ''Anio de autorizacion'': ''2013'', ''Socios o Asociados'': ''Felisa Sánchez Fernández  Alonso García Hernández esteban garcia Antonio sepulveda'', ''Rubros autorizados'': [''A. Organizaciones Civiles y Fideicomisos Asistenciales.''], ''Plantilla laboral'': ''8'', ''Plantilla voluntariado'': ''25'', ''Monto total plantilla laboral'': ''201,343'', ''Activo'': ''13,163,696'', ''Pasivo'': ''164,894'', ''Capital'': ''12,998,802'', ''Gastos administracion'': ''8,497'', ''Gastos operacion'': ''880,868'', ''Gastos representacion'': ''27,322'', ''Ingresos'': [{''Ingresos'': []}, {''OtrosIngresos'': [{''Concepto'': ''cuotas de recuperacion'', ''Monto'': ''748,750''}]}], ''Egresos'': [{''Egresos'': [{''Cargo ó Nombramiento'': ''ninguno'', ''Total de percepciones netas'': ''0''}]}, {''MontosConceptos'': [{''Monto'': ''998,814'', ''Conceptos Ejercidos'': ''asilo de ancianos medicamentos alimentación estancia defuncion''}]}, {''Donativos'': [{''Donativos'': ''0'', ''Beneficiarios'': ''ninguno''}]}], ''Actividades'': [{''Actividades'': [{''Actividad'': ''estancia asilo de ancianos'', ''Número de Beneficiarios'': ''90'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''varios''}]}]}'

# Row 6
$ws.Range("A6").Value = "'" + '2018'
$ws.Range("B6").Value = '{''Mision'': ''Mejorar la calidad de vida de los ancianos en abandono y maltrato en Guadalajara y la zona metropolitana'', ''Vision'': ''Ser una organización reconocida a nivel nacional por el excelente trato y servicio atención que se le brinda a los ancianos'', ''Pagina web'': ''http:asunciondemariaseniorshome.keyrgs.com'', ''Anio de autorizacion'': ''2013'', ''Socios o Asociados'': ''Felisa sanchez fernandez juana rodriguez contreras esteban garcia reyes maria luisa Valderrama cervantes victor Alfonso garcia hernandez'', ''Rubros autorizados'': [''A. Organizaciones Civiles y Fideicomisos Asistenciales.''], ''Plantilla laboral'': ''6'', ''Plantilla voluntariado'': ''20'', ''Monto total plantilla laboral'': ''201,061'', ''Activo'': ''12,801,693'', ''Pasivo'': ''240,483'', ''Capital'': ''12,561,210'', ''Gastos administracion'': ''1,567'', ''Gastos operacion'': ''1,214,948'', ''Gastos representacion'': ''19,735'', ''Ingresos'': [{''Ingresos'': [{''Especie '': ''60431'', ''Efectivo '': '''', ''Tipo de donante'': ''Persona Moral Privada Nacional''}, {''Especie '': '''', ''Efectivo '': ''80,763'', ''Tipo de donante'': ''Persona Moral Pública Nacional''}, {''Especie '': '''', ''Efectivo '': ''570,380'', ''Tipo de donante'': ''Persona Física Nacional''}, {''Especie '': ''60431'', ''Efectivo '': '''', ''Tipo de donante'': ''Persona Moral Privada Nacional''}, {''Especie '': '''', ''Efectivo '': ''160,000'', ''Tipo de donante'': ''Persona Moral Privada Nacional''}, {''Especie '': '''', ''Efectivo '': ''5,000'', ''Tipo de donante'': ''Persona Física Nacional''}]}, {''OtrosIngresos'': [{''Concepto'': ''ninguno'', ''Monto'': ''0''}]}], ''Egresos'': [{''Egresos'': [{''Cargo ó Nombramiento'': ''ninguno'', ''Total de percepciones netas'': ''0''}]}, {''MontosConceptos'': [{''Monto'': ''0'', ''Conceptos Ejercidos'': ''ninguno''}]}, {''Donativos'': [{''Donativos'': ''0'', ''Beneficiarios'': ''ningún particular recibo donativo''}]}], ''Actividades'': [{''Actividades'': [{''Actividad'': ''estancia asilo de ancianos'', ''Número de Beneficiarios'': ''100'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''Varios''}]}]}'

# Row 7
$ws.Range("A7").Value = "'" + '2019'
$ws.Range("B7").Value = '{''Mision'': ''Mejorar la calidad de vida de los ancianos en abandono y maltrato en Guadalajara y zona Metropolitana'', ''Vision'': ''Ser Una Organizacion Reconocida a nivel nacional por el excelente trato y servicio antencion que se le brinda a los ancianos'', ''Pagina web'': ''http://asunciondemariaseniorshome.keyrgs.com'', ''Anio de autorizacion'': ''2013'', ''Socios o Asociados'': ''Jaime Suarez Valencia, Pdte. Felisa Sanchez Fernandez, Srio. Esteban Garcia Reyes, Escrutadores Victor Alfonso Garcia Hernandez y Juana Rodriguez Contreras'', ''Rubros autorizados'': [''A. Organizaciones Civiles y Fideicomisos Asistenciales.''], ''Plantilla laboral'': ''10'', ''Plantilla voluntariado'': ''13'', ''Monto total plantilla laboral'': ''166,797'', ''Activo'': ''13,155,681'', ''Pasivo'': ''433,466'', ''Capital'': ''12,561,210'', ''Gastos administracion'': ''1,403'', ''Gastos operacion'': ''908,913'', ''Gastos representacion'': ''27,482'', ''Ingresos'': [{''Ingresos'': [{''Especie '': '''', ''Efectivo '': ''45,000'', ''Tipo de donante'': ''Persona Física Nacional''}, {''Especie '': ''Medicamentos y Despensa por 56736'', ''Efectivo '': '''', ''Tipo de donante'': ''Persona Moral Privada Nacional''}, {''Especie '': '''', ''Efectivo '': ''30,000'', ''Tipo de donante'': ''Persona Moral Privada Nacional''}, {''Especie '': '''', ''Efectivo '': ''994,215'', ''Tipo de donante'': ''Persona Física Nacional''}, {''Especie '': '''', ''Efectivo '': ''56,000'', ''Tipo de donante'': ''Persona Moral Pública Nacional''}]}, {''OtrosIngresos'': [{''Concepto'': ''Ninguno'', ''Monto'': ''0''}]}], ''Egresos'': [{''Egresos'': [{''Cargo ó Nombramiento'': ''Ninguno'', ''Total de percepciones netas'': ''0''}]}, {''MontosConceptos'': [{''Monto'': ''0'', ''Conceptos Ejercidos'': ''Ninguno''}]}, {''Donativos'': [{''Donativos'': ''0'', ''Beneficiarios'': ''Ninguno''}]}], ''Actividades'': [{''Actividades'': [{''Actividad'': ''Estancia asilo de ancianos'', ''Número de Beneficiarios'': ''120'', ''Entidad Federativa'': ''Jalisco'', ''Municipio'': ''varios''}]}]}'
